# Update countries & provincias Spain
# Refresh COVID-19 country statistics (paises.xlsx) and the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" banner in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 06:50"

# --- India (row 5): new case totals ---
$ws.Range("B5").Value = 7370468
$ws.Range("C5").Value = 4959
$ws.Range("D5").Value = 6448658
$ws.Range("E5").Value = 809649
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 112161

# --- Zimbabue overtakes Jamaica in total cases, swapping their ranking rows ---
# Row 115 becomes Zimbabue with its freshly updated figures.
$ws.Range("A115").Value = "Zimbabue"
$ws.Range("B115").Value = 8075
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 7669
$ws.Range("E115").Value = 175
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 231

# Row 116 becomes Jamaica, carrying the figures that used to sit in row 115.
$ws.Range("A116").Value = "Jamaica"
$ws.Range("B116").Value = 8067
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 3481
$ws.Range("E116").Value = 4426
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 160

# --- Tailandia (row 145): new case totals ---
$ws.Range("B145").Value = 3669
$ws.Range("C145").Value = 4
$ws.Range("D145").Value = 3467
$ws.Range("E145").Value = 143

# --- Butan (row 187): updated active/recovered counts ---
$ws.Range("D187").Value = 296
$ws.Range("E187").Value = 20
